# Update Iceland excess-mortality figures (rows 13 = Female, 39 = Male, 65 = Total)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: Female / Iceland ---
$ws.Range("D13").Value = 575
$ws.Range("P13").Value = -4.8
$ws.Range("R13").Value = -0.8
$ws.Range("U13").Value = "-4.8 (±12.1)"
$ws.Range("V13").Value = "-0.8% (±2.1%)"
$ws.Range("X13").Value = -66.90000000000001
$ws.Range("Y13").Value = 168.7
$ws.Range("Z13").Value = "-66.9(±168.7)"

# --- Row 39: Male / Iceland ---
$ws.Range("D39").Value = 441
$ws.Range("P39").Value = -35
$ws.Range("R39").Value = -7.4
$ws.Range("U39").Value = "-35.0 (±21.3)"
$ws.Range("V39").Value = "-7.4% (±3.9%)"
$ws.Range("X39").Value = -665.5
$ws.Range("Z39").Value = "-665.5(±405.0)"

# --- Row 65: Total / Iceland ---
$ws.Range("D65").Value = 1016
$ws.Range("P65").Value = -39.8
$ws.Range("R65").Value = -3.8
$ws.Range("S65").Value = 2.2
$ws.Range("U65").Value = "-39.8 (±25.5)"
$ws.Range("V65").Value = "-3.8% (±2.2%)"
$ws.Range("X65").Value = -320.2
$ws.Range("Y65").Value = 205.1
$ws.Range("Z65").Value = "-320.2(±205.1)"
